$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Name" column (B) for rows 3-9 so that each row's label
#     effectively shifts to the next surviving name (b->d, c->ffff..., d->g,
#     e->testname, ffff...->tt, g->zzzz, h->tجربة). Row 2 (a) is untouched.
$ws.Range("B3").Value = "d"
$ws.Range("B4").Value = "ffffffffffffffffffffffffffffffffffffffffffffffffffffffffffffff"
$ws.Range("B5").Value = "g"
$ws.Range("B6").Value = "testname"
$ws.Range("B7").Value = "tt"
$ws.Range("B8").Value = "zzzz"
$ws.Range("B9").Value = "تجربة"

# --- Update the price columns (C, D) for the surviving rows.
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 2.25

$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 0

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 15

$ws.Range("C7").Value = 10.025
$ws.Range("D7").Value = 15

$ws.Range("C8").Value = 1.5
$ws.Range("D8").Value = 2.01

$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 2

# --- The last five source rows (originally holding i, testname, tt, zzzz,
#     تجربة) are now finished/duplicated further up the sheet, so remove them
#     entirely.
$ws.Range("A10:D14").EntireRow.Delete()
